# Updated cryptos list on Tue May 30 17:34:49 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range("D2").Value = "27.651.20"
$ws.Range("E2").Value = "  -0.14%  "

# Row 3: update D3, E3
$ws.Range("D3").Value = "1.899.34"
$ws.Range("E3").Value = "  +0.61%  "

# Row 4: update D4, E4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5: update D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.14"
$ws.Range("E5").Value = "  -0.33%  "

# Row 6: update D6, E6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.11%  "

# Row 7: update D7, E7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5234"
$ws.Range("E7").Value = "  +8.26%  "

# Row 8: update D8, E8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3773"
$ws.Range("E8").Value = "  -0.27%  "

# Row 9: update D9, E9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07230"
$ws.Range("E9").Value = "  -1.34%  "

# Row 10: update D10, E10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.04"
$ws.Range("E10").Value = "  +3.04%  "

# Row 11: update D11, E11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8926"
$ws.Range("E11").Value = "  -2.85%  "

# Row 12: update D12, E12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07617"
$ws.Range("E12").Value = "  -0.67%  "

# Row 13: update D13, E13
$ws.Range("D13").Value = "1.892.59"
$ws.Range("E13").Value = "  +0.15%  "

# Row 14: update D14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.435"

# Row 15: update D15, E15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.83"
$ws.Range("E15").Value = "  +1.12%  "

# Row 16: update E16
$ws.Range("E16").Value = "  -0.19%  "

# Row 17: update D17, E17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008701"
$ws.Range("E17").Value = "  -1.04%  "

# Row 18: update D18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9993"

# Row 19: update D19, E19
$ws.Range("D19").Value = "27.689.45"
$ws.Range("E19").Value = "  -0.11%  "

# Row 20: update D20, E20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.43"
$ws.Range("E20").Value = "  -0.61%  "

# Row 21: update D21, E21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.122"
$ws.Range("E21").Value = "  +0.17%  "

# Row 22: update D22, E22
$ws.Range("D22").Value = "2.135.59"
$ws.Range("E22").Value = "  -0.42%  "

# Row 23: update D23, E23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.80"
$ws.Range("E23").Value = "  -0.04%  "

# Row 24: update D24, E24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.562"
$ws.Range("E24").Value = "  -0.41%  "

# Row 25: update D25, E25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.06"
$ws.Range("E25").Value = "  -0.15%  "

# Row 26: update D26, E26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.862"
$ws.Range("E26").Value = "  -2.25%  "

# Row 27: update B27, C27, D27, E27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.26"
$ws.Range("E27").Value = "  -0.44%  "

# Row 28: update B28, C28, D28, E28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.155"
$ws.Range("E28").Value = "  +2.13%  "

# Row 29: update D29, E29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.36"
$ws.Range("E29").Value = "  -1.15%  "

# Row 30: update D30, E30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.825"
$ws.Range("E30").Value = "  -1.38%  "

# Row 31: update D31, E31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08989"
$ws.Range("E31").Value = "  +0.63%  "

# Row 32: update D32, E32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.168"
$ws.Range("E32").Value = "  +0.50%  "

# Row 33: update D33, E33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.823"
$ws.Range("E33").Value = "  +4.32%  "

# Row 34: update D34, E34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.232"
$ws.Range("E34").Value = "  +1.05%  "

# Row 35: update D35, E35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7666"
$ws.Range("E35").Value = "  +0.81%  "

# Row 36: update B36, C36, D36, E36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02072"
$ws.Range("E36").Value = "  +1.61%  "

# Row 37: update B37, C37, D37, E37
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.595"
$ws.Range("E37").Value = "  +2.34%  "

# Row 38: update D38, E38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.054"
$ws.Range("E38").Value = "  +2.71%  "

# Row 39: update E39
$ws.Range("E39").Value = "  -0.11%  "

# Row 40: update D40, E40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5475"
$ws.Range("E40").Value = "  +0.44%  "

# Row 41: update D41, E41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.05262"
$ws.Range("E41").Value = "  +0.38%  "

# Row 42: update D42, E42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.622"
$ws.Range("E42").Value = "  -4.59%  "

# Row 43: update D43, E43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "113.04"
$ws.Range("E43").Value = "  +3.36%  "

# Row 44: update D44, E44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.425"
$ws.Range("E44").Value = "  +1.45%  "

# Row 45: update E45
$ws.Range("E45").Value = "  -0.84%  "

# Row 46: update D46, E46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4764"
$ws.Range("E46").Value = "  -0.16%  "

# Row 47: update D47, E47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.37"
$ws.Range("E47").Value = "  -1.92%  "

# Row 48: update D48, E48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9994"
$ws.Range("E48").Value = "  -0.16%  "

# Row 49: update D49, E49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.608"
$ws.Range("E49").Value = "  -1.11%  "

# Row 50: update D50, E50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.30"
$ws.Range("E50").Value = "  -1.67%  "

# Row 51: update D51, E51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05996"
$ws.Range("E51").Value = "  -0.93%  "

